$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 107, pushing the existing data (rows 107-190) down to 108-191.
$ws.Rows("107:107").Insert()

# Populate the newly inserted row 107 with this week's new record.
$ws.Range("A107").Value2 = 6
$ws.Range("B107").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C107").Value2 = "Metropolitana"
$ws.Range("D107").Value2 = 44634
$ws.Range("E107").Value2 = 13
$ws.Range("F107").Value2 = 100112001
$ws.Range("G107").Value2 = "Berenjena"
$ws.Range("H107").Value2 = "Sin especificar"
$ws.Range("I107").Value2 = "Primera"
$ws.Range("J107").Value2 = 340
$ws.Range("K107").Value2 = 7000
$ws.Range("L107").Value2 = 8000
$ws.Range("M107").Value2 = 7529
$ws.Range("N107").Value2 = "$/caja 50 unidades"
$ws.Range("O107").Value2 = "Región de Arica y Parinacota"
$ws.Range("P107").Value2 = 151
$ws.Range("Q107").Value2 = 50
$ws.Range("R107").Value2 = "Hortaliza"
